$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.663.34'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.295.39'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +19.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.10'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.617'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.96%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.619'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.22'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +5.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0940'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.66'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +10.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.62'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.22%  '
$ws.Range('D15').Value = '2.634.51'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.850'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '2.290.09'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '43.656.37'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.58'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.54'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.48'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.56'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.86'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.61'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.32'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.71%  '
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.30'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0933'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.64'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.55'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.127'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.72'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.80%  '
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0358'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.84'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +13.12%  '
$ws.Range('B40').Value = 'MultiversX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '74.15'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +14.95%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.42'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.88%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.83'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +12.82%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.243'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.46'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +8.45%  '
$ws.Range('E45').Value = '  +14.83%  '
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.78'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.93'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.101'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.31%  '
$ws.Range('E50').Value = '  +4.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.453'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.82%  '
